# ModelApproachSensorTimeSeriesData.pptx - "minor update to include 1D CNN eng features"
#
# 1) Refresh the cached "datetimeFigureOut" date placeholder text (slide master +
#    every slide layout) from 2/18/2020 -> 2/25/2020.
# 2) On slide 1, grow the bulleted summary textbox (it has spAutoFit, so its
#    height/position need to be nudged to fit the extra bullet line) and add a
#    new sub-bullet "Fast Fourier Transform + Discrete Wavelet Transform + 1D CNN"
#    right after the existing "Fast Fourier Transform + Discrete Wavelet
#    Transform" sub-bullet.

$p = $ppt.ActivePresentation

# --- 1) Date placeholders -------------------------------------------------

$newDate = "2/25/2020"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $lshp = $layout.Shapes.Item($si)
        if ($lshp.Name -like "Date Placeholder*") {
            $lshp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Slide 1 summary textbox ------------------------------------------

$slide1 = $p.Slides.Item(1)

$summaryBox = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $cand = $slide1.Shapes.Item($i)
    if ($cand.Name -eq "TextBox 24") {
        $summaryBox = $cand
    }
}

$tr = $summaryBox.TextFrame.TextRange

# Find the "Fast Fourier Transform + Discrete Wavelet Transform" bullet and
# insert the new "... + 1D CNN" bullet right after it (as its own paragraph).
$targetIdx = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $paraText = $tr.Paragraphs($i).Text.TrimEnd([char]13)
    if ($paraText -eq "Fast Fourier Transform + Discrete Wavelet Transform") {
        $targetIdx = $i
    }
}

$targetPara = $tr.Paragraphs($targetIdx)
$targetPara.InsertAfter([char]13 + "Fast Fourier Transform + Discrete Wavelet Transform + 1D CNN") | Out-Null

# Resize/reposition the textbox to its new spAutoFit extent.
$summaryBox.Left = 8.266225
$summaryBox.Top = 62.33
$summaryBox.Width = 647.5072
$summaryBox.Height = 540.42657
